# Apply odds/score updates to Sheet1 as described in the commit diff.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Row 7
$ws.Range("G7").Value = 5.5
$ws.Range("H7").Value = 5
$ws.Range("I7").Value = 1.44
$ws.Range("J7").Value = 5
$ws.Range("K7").Value = 2.75
$ws.Range("L7").Value = 1.83
$ws.Range("M7").Value = 26
$ws.Range("Q7").Value = 1.33
$ws.Range("R7").Value = 3.25
$ws.Range("Y7").Value = 19
$ws.Range("AB7").Value = 34
$ws.Range("AC7").Value = 26
$ws.Range("AD7").Value = 11
$ws.Range("AE7").Value = 15
$ws.Range("AK7").Value = 12
$ws.Range("AN7").Value = 8
$ws.Range("AO7").Value = 26
$ws.Range("AX7").Value = 7
$ws.Range("AZ7").Value = 17
# Row 8
$ws.Range("G8").Value = 2.75
$ws.Range("H8").Value = 3.2
$ws.Range("I8").Value = 2.55
$ws.Range("J8").Value = 3.25
$ws.Range("K8").Value = 2.05
$ws.Range("L8").Value = 3.1
$ws.Range("M8").Value = 1.04
$ws.Range("N8").Value = 9
$ws.Range("O8").Value = 1.33
$ws.Range("P8").Value = 3.25
$ws.Range("Q8").Value = 2.05
$ws.Range("R8").Value = 1.75
$ws.Range("S8").Value = 1.44
$ws.Range("T8").Value = 2.63
$ws.Range("X8").Value = 13
$ws.Range("Y8").Value = 11
$ws.Range("Z8").Value = 26
$ws.Range("AA8").Value = 23
$ws.Range("AB8").Value = 34
$ws.Range("AC8").Value = 9
$ws.Range("AD8").Value = 6
$ws.Range("AH8").Value = 8.5
$ws.Range("AI8").Value = 12
$ws.Range("AJ8").Value = 10
$ws.Range("AK8").Value = 23
$ws.Range("AN8").Value = 4.75
$ws.Range("AP8").Value = 26
$ws.Range("AQ8").Value = 51
$ws.Range("AT8").Value = 2.63
$ws.Range("AW8").Value = 4.5
$ws.Range("AZ8").Value = 41
# Row 11
$ws.Range("J11").Value = 2.75
$ws.Range("M11").Value = 1.13
$ws.Range("N11").Value = 6
$ws.Range("AI11").Value = 21
